# Increment the "Förändrad" (Changed) date in column C by one day
# for every data row (rows 2-145), e.g. 46075 -> 46076 (2026-02-22 -> 2026-02-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 145
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}
